# Apply edits to "best-feasible-slns" sheet:
#  - update the sheet view (scroll position / selection / zoom)
#  - fill in missing C152/C153 values (gap column D recalculates automatically)
#  - update C212:C241 values (gap column D recalculates automatically)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the visible window / selection of the sheet ---
$ws.Application.ActiveWindow.ScrollRow = 145
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Application.ActiveWindow.Zoom = 62
$ws.Range("C153").Select()

# --- New C values for rows 152-153 (previously empty) ---
$ws.Range("C152").Value = 116981
$ws.Range("C153").Value = 118857

# --- Updated C values for rows 212-241 ---
$newValues = @{
    212 = 56374
    213 = 57962
    214 = 56291
    215 = 56563
    216 = 56279
    217 = 56814
    218 = 55769
    219 = 56079
    220 = 56873
    221 = 55951
    222 = 107247
    223 = 107679
    224 = 105989
    225 = 106215
    226 = 106880
    227 = 106726
    228 = 105880
    229 = 103667
    230 = 106280
    231 = 105409
    232 = 149885
    233 = 149616
    234 = 152610
    235 = 152826
    236 = 150115
    237 = 148038
    238 = 147215
    239 = 152591
    240 = 149351
    241 = 149294
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 3).Value = $newValues[$row]
}

$ws.Range("C153").Select()
